# Update the "K" column (column G) values for data rows 2-20.
# These raw values were regenerated upstream (Strike# -> K recalculation,
# along with std/mean + s_vals recompute), so here we simply overwrite the
# stored numbers with the newly computed results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 2
    4  = 6
    5  = 1
    6  = 7
    7  = 1
    8  = 2
    9  = 1
    10 = 2
    11 = 0
    12 = 2
    13 = 2
    14 = 2
    15 = 3
    16 = 3
    17 = 0
    18 = 3
    19 = 3
    20 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
